$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Data corrections for a few recent rows (source data revised)
# ---------------------------------------------------------------------
$ws.Range("B112").Value2 = 3935
$ws.Range("B115").Value2 = 4111
$ws.Range("B117").Value2 = 3918

# ---------------------------------------------------------------------
# 2) Insert a new row before the last row (118) - Excel shifts the old
#    row 118 down to row 119 (values/formulas/style move with it, with
#    references auto-adjusted).
# ---------------------------------------------------------------------
$ws.Rows.Item(118).Insert()

# ---------------------------------------------------------------------
# 3) Copy the formatting (number formats / fill / style) from row 117
#    down onto the newly inserted (blank) row 118, since it is no longer
#    the "latest" row.
# ---------------------------------------------------------------------
$ws.Range("A117:M117").Copy()
$ws.Range("A118:M118").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Rows.Item(118).RowHeight = 14

# ---------------------------------------------------------------------
# 4) Fill in row 118 with the data that used to belong to the old last
#    row (now shifted to 119's position conceptually) -- i.e. restate
#    the same Date/Doses that day 44303 had.
# ---------------------------------------------------------------------
$ws.Range("A118").Value2 = 44303
$ws.Range("B118").Value2 = 3956
$ws.Range("C118").Formula = "=(AVERAGE(B112:B118))"
$ws.Range("D118").Formula = "=(D117-B118)"
$ws.Range("E118").Formula = "=E117+B118"
$ws.Range("F118").Formula = "=(E118-G118)"
$ws.Range("G118").Value2 = 14306
$ws.Range("H118").Formula = "=(G118*2)"
$ws.Range("I118").Formula = "=G118/2"
$ws.Range("J118").Formula = "=D118/C118"
$ws.Range("K118").Formula = "=A118+J118"
$ws.Range("L118").Formula = "=D118/84"
$ws.Range("M118").Value2 = ""

# ---------------------------------------------------------------------
# 5) Row 117's D formula simplifies (no longer subtracts I117) and I117
#    becomes a live formula instead of a literal 0.
# ---------------------------------------------------------------------
$ws.Range("D117").Formula = "=(D116-B117)"
$ws.Range("I117").Formula = "=G117/2"

# ---------------------------------------------------------------------
# 6) New last row (119) with the new day's data - inherits the "Good"
#    (green) styling that the old row 118 used to carry (it moved down
#    with the Insert above), so only content needs populating.
# ---------------------------------------------------------------------
$ws.Range("A119").Value2 = 44304
$ws.Range("B119").Value2 = 3541
$ws.Range("C119").Formula = "=(AVERAGE(B113:B119))"
$ws.Range("D119").Formula = "=(D118-B119)"
$ws.Range("E119").Formula = "=E118+B119"
$ws.Range("F119").Formula = "=(E119-G119)"
$ws.Range("G119").Value2 = 14306
$ws.Range("H119").Formula = "=(G119*2)"
$ws.Range("I119").Formula = "=G119/2"
$ws.Range("J119").Formula = "=D119/C119"
$ws.Range("K119").Formula = "=A119+J119"
$ws.Range("L119").Formula = "=D119/84"
$ws.Range("M119").Value2 = "daily rate to achieve June 20 target"

# ---------------------------------------------------------------------
# 7) View state: selection moved to E129.
# ---------------------------------------------------------------------
$ws.Range("E129").Select()

$wb.Application.Calculate()
